$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (F1:H1), bold like the other headers ---
$ws.Range("F1").Value = "Optimistic (Hours)"
$ws.Range("G1").Value = "Most Likely (Hours)"
$ws.Range("H1").Value = "Pessimistic (Hours)"
$ws.Range("F1:H1").Font.Bold = $true

# --- Row 2 gets its own (non-shared) formulas, exactly like Excel would
#     produce when the first formula in a column is entered on its own ---
$ws.Range("F2").Formula = "=(B2*0.5)"
$ws.Range("H2").Formula = "=(B2*1.5)"

# --- Rows 3:49 become one shared-formula block per column when filled
#     together, matching the ref="F3:F49"/ref="H3:H49" shared groups ---
$ws.Range("F3:F49").Formula = "=(B3*0.5)"
$ws.Range("H3:H49").Formula = "=(B3*1.5)"

# --- Column G ("Most Likely (Hours)") holds literal values equal to the
#     Story Point value (B) for each row -- i.e. 2 * Optimistic ---
$gValues = @(5, 5, 2, 2, 3, 3, 3, 2, 1, 1, 3, 5, 3, 5, 5, 2, 3, 3, 3, 1, 5, 2, 5, 3, 8, 5, 5, 8, 5, 2, 8, 3, 3, 5, 3, 3, 8, 5, 8, 5, 5, 3, 8, 5, 8, 5, 5, 5)
$r = 2
foreach ($v in $gValues) {
    $ws.Cells.Item($r, 7).Value = $v
    $r = $r + 1
}

# --- Widen the new columns (and E, whose header text now sits next to the
#     new table) to fit their header text, mirroring Excel's autofit ---
$ws.Columns.Item(5).ColumnWidth = 20.75
$ws.Columns.Item(6).ColumnWidth = 16.59
$ws.Columns.Item(7).ColumnWidth = 17.75
$ws.Columns.Item(8).ColumnWidth = 16.92

# --- Leave the selection where Excel would after tabbing off the filled
#     range (one column past the last entry, back on row 1) ---
$ws.Range("I1").Select()
